$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2419.3
$ws.Range("C2").Value = 216

$ws.Range("B3").Value = 289
$ws.Range("C3").Value = 216

$ws.Range("B5").Value = 1469
$ws.Range("C5").Value = 216

$ws.Range("B6").Value = 1148
$ws.Range("C6").Value = 216

$ws.Range("B7").Value = 390
$ws.Range("C7").Value = 216

$ws.Range("B9").Value = 1339
$ws.Range("C9").Value = 216

$ws.Range("B10").Value = 202
$ws.Range("C10").Value = 108

$ws.Range("B11").Value = 6
$ws.Range("C11").Value = 108

$ws.Range("B12").Value = 157
$ws.Range("C12").Value = 216
